# Reenable macroeconomic feedback loops:
# The BDMFL sheet's B2 cell holds the boolean control lever value.
# 0 = macroeconomic feedback loops are enabled (default behavior)
# 1 = macroeconomic feedback loops are disabled
# Setting it back to 0 re-enables the feedback loops.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BDMFL")
$ws.Range("B2").Value = 0
